$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Id" column header (request version of data download - instrument id)
$ws.Range("I1").Value = "Id"

# Per-instrument id values returned by the download, one per row
$ws.Range("I2").Value = 3564
$ws.Range("I3").Value = 528
$ws.Range("I4").Value = 518
$ws.Range("I5").Value = 511
$ws.Range("I6").Value = 608
$ws.Range("I7").Value = 540
$ws.Range("I8").Value = 546
$ws.Range("I9").Value = 604
$ws.Range("I10").Value = 3263
$ws.Range("I11").Value = 11141
$ws.Range("I12").Value = 639
$ws.Range("I13").Value = 703
$ws.Range("I14").Value = 702
$ws.Range("I15").Value = 4664
$ws.Range("I16").Value = 748
$ws.Range("I17").Value = 12625
$ws.Range("I18").Value = 6005
$ws.Range("I19").Value = 477
$ws.Range("I20").Value = 6461
$ws.Range("I21").Value = 6907
$ws.Range("I22").Value = 6494
$ws.Range("I23").Value = 4042
$ws.Range("I24").Value = 6329
$ws.Range("I25").Value = 8208
$ws.Range("I26").Value = 424
$ws.Range("I27").Value = 9864
$ws.Range("I28").Value = 4338
$ws.Range("I29").Value = 7217
$ws.Range("I30").Value = 6515

# Stray marker left over from the manual lookup of row 24 while testing
$ws.Range("K24").Value = " "

# Leave the selection where editing finished
$ws.Range("I22").Select()
